$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update EAN codes in A2 and A3 (leading apostrophe forces text storage,
# matching the original inlineStr/text cell type); reset to the Normal
# style afterward so no quote-prefix formatting is introduced.
$ws.Range("A2").Value = "'4517820972430"
$ws.Range("A3").Value = "'4892958240271"
$ws.Range("A2:A3").Style = "Normal"

# Remove the now-unused rows 4 and 5 (values only existed there previously)
$ws.Range("A4:A5").EntireRow.Delete()
